$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct Ben's birthday from 12/22/2020 to 12/21/2020
$ws.Range("C10").Value = "12/21/2020"

# Update the selected/active cell in the sheet view
$ws.Range("E14").Select()
